$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from the existing H1 header cell onto the new ones
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data cells (rows 2 and 3)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 6
